$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe the old data region (A1:F13) completely: content + formatting ---
$ws.Range("A1:F13").Clear()

# --- Header row (B1:E1): Input, Parameters, Loss, Accuracy, all centered ---
$ws.Range("B1").Value = "Input"
$ws.Range("C1").Value = "Parameters"
$ws.Range("D1").Value = "Loss"
$ws.Range("E1").Value = "Accuracy"

# --- Data rows: Model, Input, Parameters, Loss, Accuracy ---
$data = @(
    @("MobileNet(alpha=0.25)",   224, "0.2M", 0.75434976816177302, 0.71465969085693304),
    @("MobileNetV2(alpha=0.35)", 224, "0.4M", 0.34409651160240101, 0.88132637739181496),
    @("MobileNetV2(alpha=0.50)", 224, "0.7M", 0.30661496520042397, 0.89965093135833696),
    @("MobileNet(alpha=0.50)",   224, "0.8M", 0.54651391506195002, 0.791448533535003),
    @("MobileNetV3(small)",      224, "0.9M", 0.29235821962356501, 0.90663176774978604),
    @("MobileNetV2(alpha=0.75)", 224, "1.4M", 0.30794841051101601, 0.89790576696395796),
    @("MobileNet(alpha=0.75)",   224, "1.8M", 0.52162706851959195, 0.80715531110763505),
    @("MobileNetV2(alpha=1.0)",  224, "2.3M", 0.35025388002395602, 0.89267015457153298),
    @("MobileNetV3(large)",      224, "3.0M", 0.25453931093215898, 0.93019199371337802),
    @("MobileNet(alpha=1.0)",    224, "3.2M", 0.46997523307800199, 0.83595114946365301)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# --- Number formats: Loss/Accuracy columns (D,E) show 2 decimals ---
$ws.Range("D2:E11").NumberFormat = "0.00_ "

# --- Alignment: everything in the table is centered (except column A, left as default) ---
$ws.Range("B1:E1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B2:E11").HorizontalAlignment = -4108  # xlCenter

# --- Trailing blank style-only row (row 13): left-aligned placeholder style ---
$ws.Range("A13").HorizontalAlignment = -4131     # xlLeft
$ws.Range("D13").HorizontalAlignment = -4131     # xlLeft
$ws.Range("E13").HorizontalAlignment = -4131     # xlLeft

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 25
$ws.Columns.Item(2).ColumnWidth = 6
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 5.5
$ws.Columns.Item(5).ColumnWidth = 9

# --- Selection back to A1 (no special selection marker) ---
$ws.Range("A1").Select()
